# ---------------------------------------------------------------------------
# Renames the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# (commit: "Use <formatversion> as suffix for table headers"), then turns the
# header range A1:U57 into a proper Excel Table (ListObject) and freezes the
# header row, matching the regenerated AHB-diff workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Rename the column headers: "<Name>_old" -> "<Name>_FV2210"
#    and "<Name>_new" -> "<Name>_FV2304". Column K ("diff") is left as-is.
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value2 = $baseNames[$i] + "_FV2210"
}

$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value2 = $baseNames[$i] + "_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Convert A1:U57 into an Excel Table ("Table1"), keeping the existing
#    header-row look (bold / shaded / centred / wrapped / bordered) instead
#    of letting the new table steal it into a header dxf override. We do
#    this by stashing the current header formatting on a scratch range,
#    clearing the header range (so ListObjects.Add sees "no formatting" and
#    doesn't synthesize a dxf), adding the table, then restoring the look.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$headerRange.Copy() | Out-Null
$scratchRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$headerRange.ClearFormats()

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)

$scratchRange.Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$scratchRange.Clear()

$lo.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
